$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update a few values in rows 2-3
$ws.Range("B2").Value = 108
$ws.Range("G2").Value = 70
$ws.Range("L2").Value = 0

$ws.Range("B3").Value = 203
$ws.Range("G3").Value = 196
$ws.Range("L3").Value = 192

# Clear the contents (keep formatting) of rows 4-11 for columns B, G, L
$ws.Range("B4:B11").ClearContents()
$ws.Range("G4:G11").ClearContents()
$ws.Range("L4:L11").ClearContents()

# Update the active selection to B3
$ws.Range("B3").Select()
